# Resize the columns of the "SoundManager" calling-points table
# (the first table in the document): narrower "Calling Point" /
# wider "Function called" & "Notes" columns, and switch the table
# to a fixed layout so the widths stick.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Lock the column widths instead of letting Word auto-fit them.
$t.AllowAutoFit = $false

# New widths (dxa -> points, 1 pt = 20 dxa):
#   col 1 "Function called": 3976 -> 3438 dxa (171.9 pt)
#   col 2 "Calling Point"  : 2790 -> 2790 dxa (139.5 pt, unchanged)
#   col 3 "Notes"          : 2319 -> 2857 dxa (142.85 pt)
$t.Columns(1).Width = 171.9
$t.Columns(2).Width = 139.5
$t.Columns(3).Width = 142.85
